# Fruta / hortaliza, semanal
# Insert a new weekly record at the top of the data block (row 224),
# shifting all existing rows (224-264) down by one (to 225-265).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 224, pushing rows 224:264 down to 225:265.
$ws.Rows.Item(224).Insert()

# Populate the newly inserted row 224 with this week's data.
$ws.Range("A224").Value = 8
$ws.Range("B224").Value = "Terminal La Palmera de La Serena"
$ws.Range("C224").Value = "Coquimbo"
$ws.Range("D224").Value = 45244
$ws.Range("E224").Value = 4
$ws.Range("F224").Value = 100112044
$ws.Range("G224").Value = "Perejil"
$ws.Range("H224").Value = "Sin especificar"
$ws.Range("I224").Value = "Primera"
$ws.Range("J224").Value = 2000
$ws.Range("K224").Value = 2000
$ws.Range("L224").Value = 2500
$ws.Range("M224").Value = 2250
$ws.Range("N224").Value = '$/atado 1 a 1,5 kilos'
$ws.Range("O224").Value = 'Provincia del Elquí'
$ws.Range("P224").Value = 1500
$ws.Range("Q224").Value = 1.5
$ws.Range("R224").Value = "Hortaliza"
